{"js": "// Apply the tutorial-text revisions described by the commit.\n// 1. Title: \"Traducci\u00f3n de un Documento de la Base de Datos de Mainumby\"\n//           -> \"Mainumby: Traducir un Documento de la Base de Datos\"\n// 2. Remove the old \"_GoBack\" bookmark that sat right after \"se presenta\".\n// 3. \"Se presentar\u00e1 un men\u00fa de categor\u00edas.\" -> \"Se presenta un men\u00fa de categor\u00edas.\"\n// 4. \"Seleccion\u00e1 un documento por su nombre.\"\n//    -> \"Se presenta una lista de nombres de documentos. Seleccion\u00e1 uno.\"\n// 5. Insert a new sentence in front of \"Seleccion\u00e1 una oraci\u00f3n para traducir. \"\n// 6. \" documentos subidos.\" -> \" documentos subidos (ya descrita en otro tutorial).\"\n//    and re-insert the \"_GoBack\" bookmark right after \"...ya descrita\".\n\nconst body = context.document.body;\n\n// --- 1. Title -------------------------------------------------------\nconst titleResults = body.search(\n  \"Traducci\u00f3n de un Documento de la Base de Datos de Mainumby\",\n  { matchCase: true }\n);\ntitleResults.load(\"items\");\nawait context.sync();\ntitleResults.items[0].insertText(\n  \"Mainumby: Traducir un Documento de la Base de Datos\",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- 2. Drop the original \"_GoBack\" bookmark -------------------------\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// --- 3. \"presentar\u00e1\" -> \"presenta\" -----------------------------------\nconst menuResults = body.search(\n  \"presentar\u00e1 un men\u00fa de categor\u00eda\",\n  { matchCase: true }\n);\nmenuResults.load(\"items\");\nawait context.sync();\nmenuResults.items[0].insertText(\"presenta un men\u00fa de categor\u00eda\", \"Replace\");\nawait context.sync();\n\n// --- 4. \"Seleccion\u00e1 un documento por su nombre.\" ----------------------\nconst pickResults = body.search(\n  \"Seleccion\u00e1 un documento por su nombre.\",\n  { matchCase: true }\n);\npickResults.load(\"items\");\nawait context.sync();\npickResults.items[0].insertText(\n  \"Se presenta una lista de nombres de documentos. Seleccion\u00e1 uno.\",\n  \"Replace\"\n);\nawait context.sync();\n\n// --- 5. Prepend a new sentence before \"Seleccion\u00e1 una oraci\u00f3n...\" ----\nconst sentenceResults = body.search(\n  \"Seleccion\u00e1 una oraci\u00f3n para traducir. \",\n  { matchCase: true }\n);\nsentenceResults.load(\"items\");\nawait context.sync();\nsentenceResults.items[0].insertText(\n  \"Se muestra el contenido del documento en el espacio a la izquierda, segmentado en oraciones. \",\n  \"Before\"\n);\nawait context.sync();\n\n// --- 6. \" documentos subidos.\" -> \"... (ya descrita en otro tutorial).\"\nconst uploadedResults = body.search(\n  \" documentos subidos.\",\n  { matchCase: true }\n);\nuploadedResults.load(\"items\");\nawait context.sync();\nuploadedResults.items[0].insertText(\n  \" documentos subidos (ya descrita en otro tutorial).\",\n  \"Replace\"\n);\nawait context.sync();\n\n// Re-anchor \"_GoBack\" right after \"...ya descrita\".\nconst describedResults = body.search(\"ya descrita\", { matchCase: true });\ndescribedResults.load(\"items\");\nawait context.sync();\nconst describedEnd = describedResults.items[0].getRange(\"End\");\ndescribedEnd.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Apply the tutorial-text revisions described by the commit.\n# 1. Title: \"Traducci\u00f3n de un Documento de la Base de Datos de Mainumby\"\n#           -> \"Mainumby: Traducir un Documento de la Base de Datos\"\n# 2. Remove the old \"_GoBack\" bookmark that sat right after \"se presenta\".\n# 3. \"Se presentar\u00e1 un men\u00fa de categor\u00edas.\" -> \"Se presenta un men\u00fa de categor\u00edas.\"\n# 4. \"Seleccion\u00e1 un documento por su nombre.\"\n#    -> \"Se presenta una lista de nombres de documentos. Seleccion\u00e1 uno.\"\n# 5. Insert a new sentence in front of \"Seleccion\u00e1 una oraci\u00f3n para traducir. \"\n# 6. \" documentos subidos.\" -> \" documentos subidos (ya descrita en otro tutorial).\"\n#    and re-insert the \"_GoBack\" bookmark right after \"...ya descrita\".\n\n$wdReplaceAll   = 2\n$wdCollapseEnd  = 0\n\n$d = $word.ActiveDocument\n\nfunction Replace-DocText {\n    param([string]$SearchText, [string]$ReplaceText)\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Execute(\n        $SearchText, $false, $false, $false, $false, $false, $true, 1, $false,\n        $ReplaceText, $wdReplaceAll\n    ) | Out-Null\n}\n\n# --- 1. Title ----------------------------------------------------------\nReplace-DocText `\n    \"Traducci\u00f3n de un Documento de la Base de Datos de Mainumby\" `\n    \"Mainumby: Traducir un Documento de la Base de Datos\"\n\n# --- 2. Drop the original \"_GoBack\" bookmark ----------------------------\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# --- 3. \"presentar\u00e1\" -> \"presenta\" --------------------------------------\nReplace-DocText `\n    \"presentar\u00e1 un men\u00fa de categor\u00eda\" `\n    \"presenta un men\u00fa de categor\u00eda\"\n\n# --- 4. \"Seleccion\u00e1 un documento por su nombre.\" ------------------------\nReplace-DocText `\n    \"Seleccion\u00e1 un documento por su nombre.\" `\n    \"Se presenta una lista de nombres de documentos. Seleccion\u00e1 uno.\"\n\n# --- 5. Prepend a new sentence before \"Seleccion\u00e1 una oraci\u00f3n...\" ------\nReplace-DocText `\n    \"Seleccion\u00e1 una oraci\u00f3n para traducir. \" `\n    \"Se muestra el contenido del documento en el espacio a la izquierda, segmentado en oraciones. Seleccion\u00e1 una oraci\u00f3n para traducir. \"\n\n# --- 6. \" documentos subidos.\" -> \"... (ya descrita en otro tutorial).\" -\nReplace-DocText `\n    \" documentos subidos.\" `\n    \" documentos subidos (ya descrita en otro tutorial).\"\n\n# Re-anchor \"_GoBack\" right after \"...ya descrita\".\n$rng2 = $d.Content\n$find2 = $rng2.Find\n$find2.ClearFormatting()\n$find2.Execute(\"ya descrita\") | Out-Null\n$rng2.Collapse($wdCollapseEnd)\n$d.Bookmarks.Add(\"_GoBack\", $rng2)\n"}
